$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in rows 2-4
$ws.Range("B2").Value = 444

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 436

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 122

# Remove row 5 entirely (was A5=2, B5=183)
$ws.Range("A5:B5").Delete()
